$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("payment-request")

# The form block (header title + field-name row + sample-value row + type row)
# is duplicated three times on the "payment-request" sheet, one blank
# separator row apart, to hold a list of sheet-query entries. A new
# "insertion" label is also added next to the title in each block.

$rowOffsets = @(0, 5, 10)

foreach ($offset in $rowOffsets) {
    $r1 = 1 + $offset
    $r2 = 2 + $offset
    $r3 = 3 + $offset
    $r4 = 4 + $offset

    # Row 1: title + new "insertion" label
    $ws.Cells.Item($r1, 1).Value = "dobpaymentrequest"
    $ws.Cells.Item($r1, 2).Value = "insertion"

    # Row 2: field names
    $ws.Cells.Item($r2, 1).Value = "id"
    $ws.Cells.Item($r2, 2).Value = "code"
    $ws.Cells.Item($r2, 3).Value = "creationdate"
    $ws.Cells.Item($r2, 4).Value = "modifieddate"
    $ws.Cells.Item($r2, 5).Value = "creationinfo"
    $ws.Cells.Item($r2, 6).Value = "modificationinfo"
    $ws.Cells.Item($r2, 7).Value = "currentstates"
    $ws.Cells.Item($r2, 8).Value = "paymentType"
    $ws.Cells.Item($r2, 9).Value = "purchaseUnitId"
    $ws.Cells.Item($r2, 10).Value = "paymentForm"
    $ws.Cells.Item($r2, 11).Value = "amount"

    # Row 3: sample values. A leading apostrophe is Excel's own "force text"
    # input marker and gets consumed rather than stored, so these literal
    # values (which must keep their leading quote character) are entered
    # with it doubled up - matching how Excel itself treats a doubled
    # leading quote as an escaped, literal one.
    $ws.Cells.Item($r3, 1).Value = 1001
    $ws.Cells.Item($r3, 2).Value = "''2019000001'"
    $ws.Cells.Item($r3, 3).Value = "''2018-08-05 09:02:00'"
    $ws.Cells.Item($r3, 4).Value = "''2018-08-05 09:02:00'"
    $ws.Cells.Item($r3, 5).Value = "''Amira.Atya'"
    $ws.Cells.Item($r3, 6).Value = "''Amira.Atya'"
    $ws.Cells.Item($r3, 7).Value = "''[""Draft""]'"
    $ws.Cells.Item($r3, 8).Value = "''GENERAL'"
    $ws.Cells.Item($r3, 9).Value = 20
    $ws.Cells.Item($r3, 10).Value = "''CASH'"
    $ws.Cells.Item($r3, 11).Value = 1000

    # Row 4: field types
    $ws.Cells.Item($r4, 1).Value = "Long"
    $ws.Cells.Item($r4, 2).Value = "String"
    $ws.Cells.Item($r4, 3).Value = "Date"
    $ws.Cells.Item($r4, 4).Value = "Date"
    $ws.Cells.Item($r4, 5).Value = "String"
    $ws.Cells.Item($r4, 6).Value = "String"
    $ws.Cells.Item($r4, 7).Value = "String"
    $ws.Cells.Item($r4, 8).Value = "String"
    $ws.Cells.Item($r4, 9).Value = "Long"
    $ws.Cells.Item($r4, 10).Value = "String"
    $ws.Cells.Item($r4, 11).Value = "Float"
}

# Selection moves to K8 (inside the second block) after the edit.
$ws.Range("K8").Select()
